# Generate Report for Handback
#
# - Marks the overview / per-language rows as "handed back" (was "Ready for
#   handoff").
# - Records the new handback target file + datetime for zh-cn and de-de.
# - Widens a few columns that now hold longer status / filename text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$handoverUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d86deb1eba5af066174d14af9fcbe54315243c43/e2e/d776a15b-3318-458f-8b3f-cfea753e3722.md"
$mdName      = "d776a15b-3318-458f-8b3f-cfea753e3722.md"
$handedBack  = "Handed back: in sync with en-US"

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$overview.Range("E2").Value = $handedBack
$overview.Range("F2").Value = $handedBack
$zhcn.Range("C2").Value     = $handedBack
$dede.Range("C2").Value     = $handedBack

# --- zh-cn: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $handoverUrl, "", "", $mdName)
$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = 15570276

$zhcn.Range("J2").Value = "d776a15b-3318-458f-8b3f-cfea753e3722.ee157fd339f6e56e19fd59a7c27df5bd54048765.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-11-09 10:35:37"

# --- de-de: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$dede.Hyperlinks.Add($dede.Range("I2"), $handoverUrl, "", "", $mdName)
$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = 15570276

$dede.Range("J2").Value = "d776a15b-3318-458f-8b3f-cfea753e3722.ee157fd339f6e56e19fd59a7c27df5bd54048765.de-de.xlf"
$dede.Range("K2").Value = "2016-11-09 10:35:57"

# --- Column widths: give the (now longer) status / file-name columns more room ---
$overview.Columns.Item(5).ColumnWidth = 29.1   # E: zh-cn status
$overview.Columns.Item(6).ColumnWidth = 29.1   # F: de-de status

$zhcn.Columns.Item(3).ColumnWidth  = 29.1      # C: Status
$zhcn.Columns.Item(9).ColumnWidth  = 39.1      # I: Latest Target File
$zhcn.Columns.Item(10).ColumnWidth = 39.1      # J: Latest Handback File

$dede.Columns.Item(3).ColumnWidth  = 29.1      # C: Status
$dede.Columns.Item(9).ColumnWidth  = 39.1      # I: Latest Target File
$dede.Columns.Item(10).ColumnWidth = 39.1      # J: Latest Handback File
